# 1. check the sheets / 2. check the named regions / 3. remove unwanted named
# regions / 4. insert a new first sheet ("Sheet0") and make it the selected
# (active) sheet, shifting the existing Sheet1..Sheet4 one slot to the right.
#
# The defined name "MyRange" keeps pointing at Sheet2 (workbook scope), at
# Sheet1 (sheet-local scope) and at Sheet3 (sheet-local scope); the engine
# re-keys the localSheetId automatically as sheets are inserted, so nothing
# else needs to be touched there.

$wb = $excel.ActiveWorkbook

# Insert the brand-new sheet before the current first sheet so it becomes
# sheet #1 and pushes Sheet1..Sheet4 down by one position.
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "Sheet0"

$newSheet.Range("A1").Value = "S0A1"
$newSheet.Range("B1").Value = "S0B1"

# The new sheet is the one that is selected/active; mirror that in the
# saved view state (tab selected, B1 the active cell).
$newSheet.Range("B1").Select()
$newSheet.Activate()
